$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G8").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G10").Value = 'Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G15").Value = 'nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg'
$ws.Range("G18").Value = 'Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G19").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G20").Value = 'nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G21").Value = 'yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G24").Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G25").Value = 'backup@backdoor.com, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg'
$ws.Range("G30").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G32").Value = 'Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range("G37").Value = 'nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg'
$ws.Range("G40").Value = 'Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G41").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G42").Value = 'nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G46").Value = 'nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G54").Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G62").Value = 'yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G63").Value = 'wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Range("G68").Value = 'nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G76").Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G84").Value = 'yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Range("G85").Value = 'wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Range("G90").Value = 'manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G96").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G98").Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G106").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G107").Value = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G108").Value = 'nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G112").Value = 'manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G118").Value = 'aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("G120").Value = 'Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G128").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G129").Value = 'neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G130").Value = 'nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G131").Value = 'nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G134").Value = 'asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G140").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G142").Value = 'marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G150").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G151").Value = 'wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
$ws.Range("G156").Value = 'manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg'
$ws.Range("G162").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range("G164").Value = 'marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G172").Value = 'nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G173").Value = 'wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg'
